$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNote = 'Authors report a "error probability" of 0.000. This has to be the p-value. Any p value smaller than 0.0005 rounds to 0.000. Selected this specification because Saez et al. (2009) refer to this spec in their survey of empirical ETI estimation. This is a rather high estimate. The paper by Gottfried & Schellhorn (2004) is not that conclusive as results appear to be quite unstable and sensitive to different specifications. But I am not aware of any other literature which analyses explicitely the 1990 reform.'

$ws.Range("M2").Value = $newNote
$ws.Rows.Item(2).RowHeight = 120
$ws.Range("M2").Select()
